$d = $word.ActiveDocument

# 1. Fix typo: "informatio" -> "information" in the "An added..." paragraph
$d.Content.Find.Execute("informatio with markets", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "information with markets", 2)

# 2. Fix typos: "percemt pf the market" -> "percent of the market" in the "They agreed..." paragraph
$d.Content.Find.Execute("percemt pf the market", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "percent of the market", 2)

# 3. Move the hidden "_GoBack" bookmark.
#    It currently sits alone in its own empty paragraph, right after the
#    "Representatives of embassies..." paragraph and before the "They agreed..."
#    paragraph. It needs to move into the "In the future, the effective
#    utilisation of FTAs..." paragraph, landing right after the word "of"
#    (between "...effective utilisation of" and " FTAs...").

# Locate that empty paragraph (the one holding the bookmark) via the
# surrounding text and clear it so the old bookmark goes away, while
# leaving the (now bookmark-free) empty paragraph itself intact.
$rEmpty = $d.Content
$rEmpty.Find.Execute("challenges.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$emptyPara = $rEmpty.Paragraphs(1).Next()
$emptyPara.Range.Delete()

# Find the insertion point for the new bookmark location: right after
# "the effective utilisation of" and before " FTAs".
$rTarget = $d.Content
$rTarget.Find.Execute("the effective utilisation of", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$insertPos = $rTarget.End
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
